$wb = $excel.ActiveWorkbook

# Overview sheet - Latest HO Xliff Generate Date for 7df67ebb...md row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").Value = "2016-10-24 07:07:03"

# zh-cn sheet - Correspond Handoff Datetime / Correspond Handback DateTime for row 5
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").Value = "2016-10-24 07:06:52"
$wsZhCn.Range("K5").Value = "2016-10-24 07:08:13"

# de-de sheet - Correspond Handoff Datetime / Correspond Handback DateTime for row 5
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value = "2016-10-24 07:07:03"
$wsDeDe.Range("K5").Value = "2016-10-24 07:08:29"
